$wb = $excel.ActiveWorkbook

# Rename the main data-dictionary sheet to PROSPECTIVE
$ws1 = $wb.Worksheets.Item("gearss_data_dictionary_wOMOP")
$ws1.Name = "PROSPECTIVE"

# Insert a new, empty "RETROSPECTIVE" sheet (lands right after PROSPECTIVE,
# before "notes", since "notes" is currently the active sheet)
$wsRetro = $wb.Worksheets.Add()
$wsRetro.Name = "RETROSPECTIVE"
$wsRetro.Range("F7").Select()

# Update selections on the surviving sheets to match the saved view state
$ws1.Range("A32").Select()

$wsNotes = $wb.Worksheets.Item("notes")
$wsNotes.Range("A6:H7").Select()

# Make RETROSPECTIVE the active/selected tab
$wsRetro.Select()
